$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap/rotate existing rows (B..AC) per updated odds data ---
# Row 38 <- old row 39 data
$ws.Cells.Item(38, 2).Value = 6782522
$ws.Cells.Item(38, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(38, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(38, 5).Value = 45171.75
$ws.Cells.Item(38, 6).Value = "Municipal Perez Zeledon"
$ws.Cells.Item(38, 7).Value = "Sporting San Jose"
$ws.Cells.Item(38, 8).Value = 1
$ws.Cells.Item(38, 9).Value = 2
$ws.Cells.Item(38, 10).Value = "A"
$ws.Cells.Item(38, 11).Value = 2.5
$ws.Cells.Item(38, 12).Value = 3.5
$ws.Cells.Item(38, 13).Value = 2.5
$ws.Cells.Item(38, 14).Value = 2.2
$ws.Cells.Item(38, 15).Value = 3.5
$ws.Cells.Item(38, 16).Value = 2.9
$ws.Cells.Item(38, 17).Value = -0.25
$ws.Cells.Item(38, 18).Value = 1.9
$ws.Cells.Item(38, 19).Value = 1.9
$ws.Cells.Item(38, 20).Value = 2.5
$ws.Cells.Item(38, 21).Value = 1.9
$ws.Cells.Item(38, 22).Value = 1.9
$ws.Cells.Item(38, 23).Value = -1
$ws.Cells.Item(38, 24).Value = -1
$ws.Cells.Item(38, 25).Value = 1.9
$ws.Cells.Item(38, 26).Value = -1
$ws.Cells.Item(38, 27).Value = 0.8999999999999999
$ws.Cells.Item(38, 28).Value = 0.8999999999999999
$ws.Cells.Item(38, 29).Value = -1

# Row 39 <- old row 38 data
$ws.Cells.Item(39, 2).Value = 6781354
$ws.Cells.Item(39, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(39, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(39, 5).Value = 45171.75
$ws.Cells.Item(39, 6).Value = "Puntarenas"
$ws.Cells.Item(39, 7).Value = "AD San Carlos"
$ws.Cells.Item(39, 8).Value = 1
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = "H"
$ws.Cells.Item(39, 11).Value = 2.4
$ws.Cells.Item(39, 12).Value = 3.2
$ws.Cells.Item(39, 13).Value = 2.8
$ws.Cells.Item(39, 14).Value = 2.3
$ws.Cells.Item(39, 15).Value = 3.2
$ws.Cells.Item(39, 16).Value = 3
$ws.Cells.Item(39, 17).Value = -0.25
$ws.Cells.Item(39, 18).Value = 2
$ws.Cells.Item(39, 19).Value = 1.8
$ws.Cells.Item(39, 20).Value = 2.25
$ws.Cells.Item(39, 21).Value = 1.9
$ws.Cells.Item(39, 22).Value = 1.9
$ws.Cells.Item(39, 23).Value = 1.3
$ws.Cells.Item(39, 24).Value = -1
$ws.Cells.Item(39, 25).Value = -1
$ws.Cells.Item(39, 26).Value = 1
$ws.Cells.Item(39, 27).Value = -1
$ws.Cells.Item(39, 28).Value = -1
$ws.Cells.Item(39, 29).Value = 0.8999999999999999


# Row 110 <- old row 111 data
$ws.Cells.Item(110, 2).Value = 6782579
$ws.Cells.Item(110, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(110, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(110, 5).Value = 45238.875
$ws.Cells.Item(110, 6).Value = "Santos de Gupiles"
$ws.Cells.Item(110, 7).Value = "AD San Carlos"
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 2
$ws.Cells.Item(110, 10).Value = "A"
$ws.Cells.Item(110, 11).Value = 2.4
$ws.Cells.Item(110, 12).Value = 3.3
$ws.Cells.Item(110, 13).Value = 2.7
$ws.Cells.Item(110, 14).Value = 2.375
$ws.Cells.Item(110, 15).Value = 3.4
$ws.Cells.Item(110, 16).Value = 2.8
$ws.Cells.Item(110, 17).Value = -0.25
$ws.Cells.Item(110, 18).Value = 2
$ws.Cells.Item(110, 19).Value = 1.8
$ws.Cells.Item(110, 20).Value = 2.5
$ws.Cells.Item(110, 21).Value = 1.875
$ws.Cells.Item(110, 22).Value = 1.925
$ws.Cells.Item(110, 23).Value = -1
$ws.Cells.Item(110, 24).Value = -1
$ws.Cells.Item(110, 25).Value = 1.8
$ws.Cells.Item(110, 26).Value = -1
$ws.Cells.Item(110, 27).Value = 0.8
$ws.Cells.Item(110, 28).Value = -1
$ws.Cells.Item(110, 29).Value = 0.925

# Row 111 <- old row 110 data
$ws.Cells.Item(111, 2).Value = 6782581
$ws.Cells.Item(111, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(111, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(111, 5).Value = 45238.875
$ws.Cells.Item(111, 6).Value = "Alajuelense"
$ws.Cells.Item(111, 7).Value = "AD Grecia"
$ws.Cells.Item(111, 8).Value = 2
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = "H"
$ws.Cells.Item(111, 11).Value = 1.181
$ws.Cells.Item(111, 12).Value = 6.5
$ws.Cells.Item(111, 13).Value = 11
$ws.Cells.Item(111, 14).Value = 1.25
$ws.Cells.Item(111, 15).Value = 5
$ws.Cells.Item(111, 16).Value = 9
$ws.Cells.Item(111, 17).Value = -1.75
$ws.Cells.Item(111, 18).Value = 1.975
$ws.Cells.Item(111, 19).Value = 1.825
$ws.Cells.Item(111, 20).Value = 3.25
$ws.Cells.Item(111, 21).Value = 2
$ws.Cells.Item(111, 22).Value = 1.8
$ws.Cells.Item(111, 23).Value = 0.25
$ws.Cells.Item(111, 24).Value = -1
$ws.Cells.Item(111, 25).Value = -1
$ws.Cells.Item(111, 26).Value = 0.4875
$ws.Cells.Item(111, 27).Value = -0.5
$ws.Cells.Item(111, 28).Value = -1
$ws.Cells.Item(111, 29).Value = 0.8

# Row 129 <- old row 131 data
$ws.Cells.Item(129, 2).Value = 6782596
$ws.Cells.Item(129, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(129, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(129, 5).Value = 45255.95833333334
$ws.Cells.Item(129, 6).Value = "Alajuelense"
$ws.Cells.Item(129, 7).Value = "AD Guanacasteca"
$ws.Cells.Item(129, 8).Value = 3
$ws.Cells.Item(129, 9).Value = 4
$ws.Cells.Item(129, 10).Value = "A"
$ws.Cells.Item(129, 11).Value = 1.363
$ws.Cells.Item(129, 12).Value = 4.75
$ws.Cells.Item(129, 13).Value = 8
$ws.Cells.Item(129, 14).Value = 1.444
$ws.Cells.Item(129, 15).Value = 4.333
$ws.Cells.Item(129, 16).Value = 7
$ws.Cells.Item(129, 17).Value = -1.25
$ws.Cells.Item(129, 18).Value = 1.975
$ws.Cells.Item(129, 19).Value = 1.825
$ws.Cells.Item(129, 20).Value = 2.75
$ws.Cells.Item(129, 21).Value = 1.775
$ws.Cells.Item(129, 22).Value = 2.025
$ws.Cells.Item(129, 23).Value = -1
$ws.Cells.Item(129, 24).Value = -1
$ws.Cells.Item(129, 25).Value = 6
$ws.Cells.Item(129, 26).Value = -1
$ws.Cells.Item(129, 27).Value = 0.825
$ws.Cells.Item(129, 28).Value = 0.7749999999999999
$ws.Cells.Item(129, 29).Value = -1

# Row 130 <- old row 129 data
$ws.Cells.Item(130, 2).Value = 6782595
$ws.Cells.Item(130, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(130, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(130, 5).Value = 45255.95833333334
$ws.Cells.Item(130, 6).Value = "Herediano"
$ws.Cells.Item(130, 7).Value = "Sporting San Jose"
$ws.Cells.Item(130, 8).Value = 3
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = "H"
$ws.Cells.Item(130, 11).Value = 1.4
$ws.Cells.Item(130, 12).Value = 4.75
$ws.Cells.Item(130, 13).Value = 7
$ws.Cells.Item(130, 14).Value = 1.363
$ws.Cells.Item(130, 15).Value = 4.75
$ws.Cells.Item(130, 16).Value = 8.5
$ws.Cells.Item(130, 17).Value = -1.25
$ws.Cells.Item(130, 18).Value = 1.8
$ws.Cells.Item(130, 19).Value = 2
$ws.Cells.Item(130, 20).Value = 3
$ws.Cells.Item(130, 21).Value = 1.95
$ws.Cells.Item(130, 22).Value = 1.85
$ws.Cells.Item(130, 23).Value = 0.363
$ws.Cells.Item(130, 24).Value = -1
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = 0.8
$ws.Cells.Item(130, 27).Value = -1
$ws.Cells.Item(130, 28).Value = 0
$ws.Cells.Item(130, 29).Value = -0

# Row 131 <- old row 130 data
$ws.Cells.Item(131, 2).Value = 6782598
$ws.Cells.Item(131, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(131, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(131, 5).Value = 45255.95833333334
$ws.Cells.Item(131, 6).Value = "Municipal Perez Zeledon"
$ws.Cells.Item(131, 7).Value = "Cartagines"
$ws.Cells.Item(131, 8).Value = 1
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = "H"
$ws.Cells.Item(131, 11).Value = 4.5
$ws.Cells.Item(131, 12).Value = 3.75
$ws.Cells.Item(131, 13).Value = 1.615
$ws.Cells.Item(131, 14).Value = 3.4
$ws.Cells.Item(131, 15).Value = 3.4
$ws.Cells.Item(131, 16).Value = 1.85
$ws.Cells.Item(131, 17).Value = 0.5
$ws.Cells.Item(131, 18).Value = 1.8
$ws.Cells.Item(131, 19).Value = 2
$ws.Cells.Item(131, 20).Value = 2.75
$ws.Cells.Item(131, 21).Value = 1.9
$ws.Cells.Item(131, 22).Value = 1.9
$ws.Cells.Item(131, 23).Value = 2.4
$ws.Cells.Item(131, 24).Value = -1
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = 0.8
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = -1
$ws.Cells.Item(131, 29).Value = 0.8999999999999999

# Row 200 <- old row 201 data
$ws.Cells.Item(200, 2).Value = 7624967
$ws.Cells.Item(200, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(200, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(200, 5).Value = 45353.75
$ws.Cells.Item(200, 6).Value = "Puntarenas"
$ws.Cells.Item(200, 7).Value = "Herediano"
$ws.Cells.Item(200, 8).Value = 0
$ws.Cells.Item(200, 9).Value = 0
$ws.Cells.Item(200, 10).Value = "D"
$ws.Cells.Item(200, 11).Value = 3.75
$ws.Cells.Item(200, 12).Value = 3.4
$ws.Cells.Item(200, 13).Value = 1.8
$ws.Cells.Item(200, 14).Value = 2.8
$ws.Cells.Item(200, 15).Value = 3.1
$ws.Cells.Item(200, 16).Value = 2.25
$ws.Cells.Item(200, 17).Value = 0.25
$ws.Cells.Item(200, 18).Value = 1.8
$ws.Cells.Item(200, 19).Value = 2
$ws.Cells.Item(200, 20).Value = 2.25
$ws.Cells.Item(200, 21).Value = 1.775
$ws.Cells.Item(200, 22).Value = 2.025
$ws.Cells.Item(200, 23).Value = -1
$ws.Cells.Item(200, 24).Value = 2.1
$ws.Cells.Item(200, 25).Value = -1
$ws.Cells.Item(200, 26).Value = 0.4
$ws.Cells.Item(200, 27).Value = -0.5
$ws.Cells.Item(200, 28).Value = -1
$ws.Cells.Item(200, 29).Value = 1.025

# Row 201 <- old row 200 data
$ws.Cells.Item(201, 2).Value = 7623921
$ws.Cells.Item(201, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(201, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(201, 5).Value = 45353.75
$ws.Cells.Item(201, 6).Value = "AD Grecia"
$ws.Cells.Item(201, 7).Value = "Municipal Liberia"
$ws.Cells.Item(201, 8).Value = 1
$ws.Cells.Item(201, 9).Value = 2
$ws.Cells.Item(201, 10).Value = "A"
$ws.Cells.Item(201, 11).Value = 2.75
$ws.Cells.Item(201, 12).Value = 3.25
$ws.Cells.Item(201, 13).Value = 2.3
$ws.Cells.Item(201, 14).Value = 3.1
$ws.Cells.Item(201, 15).Value = 3.25
$ws.Cells.Item(201, 16).Value = 2.1
$ws.Cells.Item(201, 17).Value = 0.25
$ws.Cells.Item(201, 18).Value = 1.9
$ws.Cells.Item(201, 19).Value = 1.9
$ws.Cells.Item(201, 20).Value = 2.5
$ws.Cells.Item(201, 21).Value = 1.9
$ws.Cells.Item(201, 22).Value = 1.9
$ws.Cells.Item(201, 23).Value = -1
$ws.Cells.Item(201, 24).Value = -1
$ws.Cells.Item(201, 25).Value = 1.1
$ws.Cells.Item(201, 26).Value = -1
$ws.Cells.Item(201, 27).Value = 0.8999999999999999
$ws.Cells.Item(201, 28).Value = 0.8999999999999999
$ws.Cells.Item(201, 29).Value = -1

# --- Add new rows 211, 212, 213 ---
# Copy formatting from row 210 (A:AC) into new rows, then set values
$ws.Range("A210:AC210").Copy()
$ws.Range("A211:AC213").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 211
$ws.Cells.Item(211, 1).Value = 209
$ws.Cells.Item(211, 2).Value = 7623995
$ws.Cells.Item(211, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(211, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(211, 5).Value = 45377.83333333334
$ws.Cells.Item(211, 6).Value = "Municipal Perez Zeledon"
$ws.Cells.Item(211, 7).Value = "Sporting San Jose"
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(211, 9).Value = 2
$ws.Cells.Item(211, 10).Value = "A"
$ws.Cells.Item(211, 11).Value = 2.3
$ws.Cells.Item(211, 12).Value = 3.2
$ws.Cells.Item(211, 13).Value = 3
$ws.Cells.Item(211, 14).Value = 3
$ws.Cells.Item(211, 15).Value = 3.2
$ws.Cells.Item(211, 16).Value = 2.4
$ws.Cells.Item(211, 17).Value = 0.25
$ws.Cells.Item(211, 18).Value = 1.775
$ws.Cells.Item(211, 19).Value = 2.025
$ws.Cells.Item(211, 20).Value = 2.5
$ws.Cells.Item(211, 21).Value = 1.95
$ws.Cells.Item(211, 22).Value = 1.85
$ws.Cells.Item(211, 23).Value = -1
$ws.Cells.Item(211, 24).Value = -1
$ws.Cells.Item(211, 25).Value = 1.4
$ws.Cells.Item(211, 26).Value = -1
$ws.Cells.Item(211, 27).Value = 1.025
$ws.Cells.Item(211, 28).Value = 0.95
$ws.Cells.Item(211, 29).Value = -1

# Row 212
$ws.Cells.Item(212, 1).Value = 210
$ws.Cells.Item(212, 2).Value = 7623993
$ws.Cells.Item(212, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(212, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(212, 5).Value = 45378.91666666666
$ws.Cells.Item(212, 6).Value = "AD San Carlos"
$ws.Cells.Item(212, 7).Value = "Herediano"
$ws.Cells.Item(212, 11).Value = 2.5
$ws.Cells.Item(212, 12).Value = 3.2
$ws.Cells.Item(212, 13).Value = 2.8
$ws.Cells.Item(212, 14).Value = 2.05
$ws.Cells.Item(212, 15).Value = 3.25
$ws.Cells.Item(212, 16).Value = 3.3
$ws.Cells.Item(212, 17).Value = -0.25
$ws.Cells.Item(212, 18).Value = 1.825
$ws.Cells.Item(212, 19).Value = 1.975
$ws.Cells.Item(212, 20).Value = 2.5
$ws.Cells.Item(212, 21).Value = 1.95
$ws.Cells.Item(212, 22).Value = 1.85
$ws.Cells.Item(212, 23).Value = 0
$ws.Cells.Item(212, 24).Value = 0
$ws.Cells.Item(212, 25).Value = 0
$ws.Cells.Item(212, 26).Value = 0
$ws.Cells.Item(212, 27).Value = 0

# Row 213
$ws.Cells.Item(213, 1).Value = 211
$ws.Cells.Item(213, 2).Value = 7623994
$ws.Cells.Item(213, 3).Value = "Costa Rica Primera Division"
$ws.Cells.Item(213, 4).Value = "Costa Rica Primera Division"
$ws.Cells.Item(213, 5).Value = 45378.95833333334
$ws.Cells.Item(213, 6).Value = "Cartagines"
$ws.Cells.Item(213, 7).Value = "Deportivo Saprissa"
$ws.Cells.Item(213, 11).Value = 3.25
$ws.Cells.Item(213, 12).Value = 3.5
$ws.Cells.Item(213, 13).Value = 2.1
$ws.Cells.Item(213, 14).Value = 3.4
$ws.Cells.Item(213, 15).Value = 3.2
$ws.Cells.Item(213, 16).Value = 2.15
$ws.Cells.Item(213, 17).Value = 0.25
$ws.Cells.Item(213, 18).Value = 1.975
$ws.Cells.Item(213, 19).Value = 1.825
$ws.Cells.Item(213, 20).Value = 2.25
$ws.Cells.Item(213, 21).Value = 1.95
$ws.Cells.Item(213, 22).Value = 1.85
$ws.Cells.Item(213, 23).Value = 0
$ws.Cells.Item(213, 24).Value = 0
$ws.Cells.Item(213, 25).Value = 0
$ws.Cells.Item(213, 26).Value = 0
$ws.Cells.Item(213, 27).Value = 0

# Rows 212 and 213 have no result yet (future matches): clear H,I,J and AB,AC
$ws.Range("H212:J213").ClearContents()
$ws.Range("AB212:AC213").ClearContents()
